$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "292.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.22%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.45%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.238"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.31%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07173"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7.94%"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.45%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.602"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.95%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.400"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.10%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9100"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.04%"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.82%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07748"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "17.46%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07708"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.49%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02920"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.78%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.26%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001612"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.16%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006549"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.21%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006079"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.30%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.487"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.27%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.76%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3252"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.12%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.45%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.043"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.37%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.41%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04532"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.83%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001207"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.29%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004262"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.62%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001166"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.68%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001683"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.01%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04438"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.69%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007015"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.92%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1279"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.43%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002202"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "9.08%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01332"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.18%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005840"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.41%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01295"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.96%"
